$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - set value then copy formatting from E1 (header style)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$timestamps = @(
    "2021-10-05 10:50:06.606765",
    "2021-10-05 10:50:06.606777",
    "2021-10-05 10:50:06.606780",
    "2021-10-05 10:50:06.606783",
    "2021-10-05 10:50:06.606786",
    "2021-10-05 10:50:06.606789",
    "2021-10-05 10:50:06.606792",
    "2021-10-05 10:50:06.606794",
    "2021-10-05 10:50:06.606797",
    "2021-10-05 10:50:06.606800",
    "2021-10-05 10:50:06.606802",
    "2021-10-05 10:50:06.606805",
    "2021-10-05 10:50:06.606808",
    "2021-10-05 10:50:06.606810",
    "2021-10-05 10:50:06.606813",
    "2021-10-05 10:50:06.606815",
    "2021-10-05 10:50:06.606818",
    "2021-10-05 10:50:06.606821",
    "2021-10-05 10:50:06.606824",
    "2021-10-05 10:50:06.606826",
    "2021-10-05 10:50:06.606829",
    "2021-10-05 10:50:06.606831",
    "2021-10-05 10:50:06.606834",
    "2021-10-05 10:50:06.606836",
    "2021-10-05 10:50:06.606839",
    "2021-10-05 10:50:06.606842",
    "2021-10-05 10:50:06.606844",
    "2021-10-05 10:50:06.606847",
    "2021-10-05 10:50:06.606850",
    "2021-10-05 10:50:06.606852",
    "2021-10-05 10:50:06.606855",
    "2021-10-05 10:50:06.606857",
    "2021-10-05 10:50:06.606860",
    "2021-10-05 10:50:06.606863"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
